$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws 'D2' '27.971.72'
Set-TextCell $ws 'E2' '  +0.38%  '
Set-TextCell $ws 'D3' '1.857.60'
Set-TextCell $ws 'E3' '  -0.31%  '
Set-TextCell $ws 'E4' '  +0.04%  '
Set-TextCell $ws 'D5' '311.68'
Set-TextCell $ws 'E5' '  -0.29%  '
Set-TextCell $ws 'D6' '1.003'
Set-TextCell $ws 'E6' '  +0.07%  '
Set-TextCell $ws 'D7' '0.5088'
Set-TextCell $ws 'E7' '  +2.11%  '
Set-TextCell $ws 'D8' '0.3806'
Set-TextCell $ws 'E8' '  -0.14%  '
Set-TextCell $ws 'D9' '0.08272'
Set-TextCell $ws 'E9' '  -6.32%  '
Set-TextCell $ws 'E10' '  -0.46%  '
Set-TextCell $ws 'B11' 'Polkadot'
Set-TextCell $ws 'C11' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws 'D11' '6.204'
Set-TextCell $ws 'E11' '  -2.27%  '
Set-TextCell $ws 'B12' 'Solana'
Set-TextCell $ws 'C12' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell $ws 'D12' '20.50'
Set-TextCell $ws 'E12' '  -0.33%  '
Set-TextCell $ws 'B13' 'WrappedEther'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws 'D13' '1.854.08'
Set-TextCell $ws 'E13' '  -0.62%  '
Set-TextCell $ws 'B14' 'Chainlink'
Set-TextCell $ws 'C14' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws 'D14' '7.200'
Set-TextCell $ws 'E14' '  -0.14%  '
Set-TextCell $ws 'B15' 'BinanceUSD'
Set-TextCell $ws 'C15' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell $ws 'D15' '1.004'
Set-TextCell $ws 'E15' '  -0.03%  '
Set-TextCell $ws 'B16' 'ShibaInu'
Set-TextCell $ws 'C16' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D16' '0.00001096'
Set-TextCell $ws 'E16' '  +0.23%  '
Set-TextCell $ws 'B17' 'Litecoin'
Set-TextCell $ws 'C17' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell $ws 'D17' '90.60'
Set-TextCell $ws 'E17' '  -0.42%  '
Set-TextCell $ws 'B18' 'TRON'
Set-TextCell $ws 'C18' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell $ws 'D18' '0.06608'
Set-TextCell $ws 'E18' '  -0.81%  '
Set-TextCell $ws 'B19' 'Avalanche'
Set-TextCell $ws 'C19' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws 'D19' '17.62'
Set-TextCell $ws 'E19' '  -1.69%  '
Set-TextCell $ws 'B20' 'Dai'
Set-TextCell $ws 'C20' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D20' '1.002'
Set-TextCell $ws 'E20' '  +0.07%  '
Set-TextCell $ws 'B21' 'Uniswap'
Set-TextCell $ws 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell $ws 'D21' '6.022'
Set-TextCell $ws 'E21' '  -1.18%  '
Set-TextCell $ws 'B22' 'WrappedBTC'
Set-TextCell $ws 'C22' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws 'D22' '27.987.85'
Set-TextCell $ws 'E22' '  +0.30%  '
Set-TextCell $ws 'B23' 'Cosmos'
Set-TextCell $ws 'C23' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D23' '11.03'
Set-TextCell $ws 'E23' '  -3.47%  '
Set-TextCell $ws 'B24' 'Toncoin'
Set-TextCell $ws 'C24' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 'D24' '2.237'
Set-TextCell $ws 'E24' '  -2.01%  '
Set-TextCell $ws 'B25' 'LidoDAOToken'
Set-TextCell $ws 'C25' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell $ws 'D25' '2.534'
Set-TextCell $ws 'E25' '  +1.40%  '
Set-TextCell $ws 'B26' 'WrappedliquidstakedEther2.0'
Set-TextCell $ws 'C26' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell $ws 'D26' '2.063.37'
Set-TextCell $ws 'E26' '  -0.72%  '
Set-TextCell $ws 'B27' 'Monero'
Set-TextCell $ws 'C27' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws 'D27' '158.08'
Set-TextCell $ws 'E27' '  +0.27%  '
Set-TextCell $ws 'B28' 'EthereumClassic'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D28' '20.44'
Set-TextCell $ws 'E28' '  -1.00%  '
Set-TextCell $ws 'B29' 'BitcoinCash'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws 'D29' '124.43'
Set-TextCell $ws 'E29' '  -1.45%  '
Set-TextCell $ws 'B30' 'Stellar'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D30' '0.1051'
Set-TextCell $ws 'E30' '  -0.49%  '
Set-TextCell $ws 'B31' 'ImmutableX'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D31' '1.034'
Set-TextCell $ws 'E31' '  -1.60%  '
Set-TextCell $ws 'B32' 'Filecoin'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D32' '5.616'
Set-TextCell $ws 'E32' '  +0.73%  '
Set-TextCell $ws 'B33' 'HuobiToken'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws 'D33' '3.592'
Set-TextCell $ws 'E33' '  -0.29%  '
Set-TextCell $ws 'B34' 'FraxShare'
Set-TextCell $ws 'C34' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws 'D34' '9.406'
Set-TextCell $ws 'E34' '  +0.16%  '
Set-TextCell $ws 'B35' 'VeChain'
Set-TextCell $ws 'C35' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D35' '0.02407'
Set-TextCell $ws 'E35' '  +0.90%  '
Set-TextCell $ws 'D36' '0.06493'
Set-TextCell $ws 'E36' '  -0.26%  '
Set-TextCell $ws 'B37' 'Algorand'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws 'D37' '0.2165'
Set-TextCell $ws 'E37' '  -0.40%  '
Set-TextCell $ws 'B38' 'ARBITRUM'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws 'D38' '1.205'
Set-TextCell $ws 'E38' '  +0.83%  '
Set-TextCell $ws 'B39' 'TheSandbox'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell $ws 'D39' '0.6431'
Set-TextCell $ws 'E39' '  +1.47%  '
Set-TextCell $ws 'B40' 'TrustWalletToken'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws 'D40' '1.223'
Set-TextCell $ws 'E40' '  -4.08%  '
Set-TextCell $ws 'B41' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D41' '4.879'
Set-TextCell $ws 'E41' '  -0.10%  '
Set-TextCell $ws 'B42' 'Aptos'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D42' '11.15'
Set-TextCell $ws 'E42' '  -2.08%  '
Set-TextCell $ws 'B43' 'Decentraland'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell $ws 'D43' '0.6070'
Set-TextCell $ws 'E43' '  +1.63%  '
Set-TextCell $ws 'B44' 'EnergySwap'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D44' '13.12'
Set-TextCell $ws 'E44' '  -0.41%  '
Set-TextCell $ws 'B45' 'WEMIXTOKEN'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws 'D45' '1.278'
Set-TextCell $ws 'E45' '  -0.30%  '
Set-TextCell $ws 'B46' 'PancakeSwap'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws 'D46' '3.650'
Set-TextCell $ws 'E46' '  -0.53%  '
Set-TextCell $ws 'B47' 'NEARProtocol'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D47' '2.008'
Set-TextCell $ws 'E47' '  +1.65%  '
Set-TextCell $ws 'B48' 'EOS'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell $ws 'D48' '1.207'
Set-TextCell $ws 'E48' '  -1.27%  '
Set-TextCell $ws 'B49' 'Quant'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell $ws 'D49' '120.00'
Set-TextCell $ws 'E49' '  -0.47%  '
Set-TextCell $ws 'B50' 'Aave'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D50' '78.57'
Set-TextCell $ws 'E50' '  +0.02%  '
Set-TextCell $ws 'B51' 'Cronos'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws 'D51' '0.06833'
Set-TextCell $ws 'E51' '  -0.80%  '
